# Add two new columns, I (I0) and J (IF), to the right of the existing
# H (IP) column. Column I is filled with a constant 1 for every data row
# and column J duplicates the value already present in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---------------------------------------------------------
# Copy H1's formatting (bold, centered, bordered) onto the new header
# cells before writing their text, same as the existing headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows ---------------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}

$excel.CutCopyMode = $false
